$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D keeps its text formatting so numeric-looking values
# (e.g. "0.999") are not auto-converted into real numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '66.623.95'
$ws.Range("E2").Value = '  +1.27%  '
$ws.Range("D3").Value = '3.230.77'
$ws.Range("E3").Value = '  +1.73%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '605.11'
$ws.Range("E5").Value = '  +1.85%  '
$ws.Range("D6").Value = '158.83'
$ws.Range("E6").Value = '  +4.21%  '
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").Value = '3.233.06'
$ws.Range("E8").Value = '  +1.97%  '
$ws.Range("D9").Value = '0.548'
$ws.Range("E9").Value = '  +2.76%  '
$ws.Range("D10").Value = '0.162'
$ws.Range("E10").Value = '  +2.01%  '
$ws.Range("D11").Value = '5.69'
$ws.Range("E11").Value = '  -5.91%  '
$ws.Range("D12").Value = '0.510'
$ws.Range("E12").Value = '  -0.15%  '
$ws.Range("D13").Value = '0.0000274'
$ws.Range("E13").Value = '  +3.24%  '
$ws.Range("D14").Value = '39.11'
$ws.Range("E14").Value = '  +1.17%  '
$ws.Range("D15").Value = '3.759.60'
$ws.Range("E15").Value = '  +1.69%  '
$ws.Range("D16").Value = '66.614.56'
$ws.Range("E16").Value = '  +1.15%  '
$ws.Range("D17").Value = '7.49'
$ws.Range("E17").Value = '  +1.26%  '
$ws.Range("D18").Value = '3.228.30'
$ws.Range("E18").Value = '  +1.78%  '
$ws.Range("E19").Value = '  +1.45%  '
$ws.Range("D20").Value = '512.41'
$ws.Range("E20").Value = '  +1.25%  '
$ws.Range("D21").Value = '15.34'
$ws.Range("E21").Value = '  +0.26%  '
$ws.Range("D22").Value = '0.737'
$ws.Range("E22").Value = '  +0.61%  '
$ws.Range("D23").Value = '8.08'
$ws.Range("E23").Value = '  +1.41%  '
$ws.Range("D24").Value = '14.84'
$ws.Range("E24").Value = '  -0.41%  '
$ws.Range("D25").Value = '84.91'
$ws.Range("E25").Value = '  +0.42%  '
$ws.Range("E26").Value = '  +0.14%  '
$ws.Range("D27").Value = '3.00'
$ws.Range("E27").Value = '  +0.91%  '
$ws.Range("D28").Value = '9.22'
$ws.Range("E28").Value = '  +0.28%  '
$ws.Range("D29").Value = '2.41'
$ws.Range("E29").Value = '  +5.99%  '
$ws.Range("E30").Value = '  +4.88%  '
$ws.Range("D31").Value = '7.06'
$ws.Range("E31").Value = '  +1.56%  '
$ws.Range("D32").Value = '28.29'
$ws.Range("E32").Value = '  +1.22%  '
$ws.Range("E33").Value = '  +0.11%  '
$ws.Range("E34").Value = '  -2.60%  '
$ws.Range("D35").Value = '6.55'
$ws.Range("E35").Value = '  +1.71%  '
$ws.Range("D36").Value = '516.65'
$ws.Range("E36").Value = '  +7.73%  '
$ws.Range("D37").Value = '0.0955'
$ws.Range("E37").Value = '  +6.32%  '
$ws.Range("D38").Value = '56.26'
$ws.Range("E38").Value = '  +2.90%  '
$ws.Range("D39").Value = '0.0₃0777'
$ws.Range("E39").Value = '  +19.88%  '
$ws.Range("D40").Value = '0.0421'
$ws.Range("E40").Value = '  +0.98%  '
$ws.Range("E41").Value = '  +7.09%  '
$ws.Range("E42").Value = '  +6.58%  '
$ws.Range("D43").Value = '8.81'
$ws.Range("E43").Value = '  +0.30%  '
$ws.Range("D44").Value = '0.301'
$ws.Range("E44").Value = '  +0.95%  '
$ws.Range("D45").Value = '2.49'
$ws.Range("E45").Value = '  +4.09%  '
$ws.Range("D46").Value = '2.872.64'
$ws.Range("E46").Value = '  -0.31%  '
$ws.Range("D47").Value = '28.67'
$ws.Range("E47").Value = '  +1.45%  '
$ws.Range("D48").Value = '2.42'
$ws.Range("D50").Value = '0.117'
$ws.Range("E50").Value = '  +0.86%  '
$ws.Range("E51").Value = '  +3.68%  '

# Restore the original (default/general) formatting & style on column D
# now that the text values are safely stored as text.
$ws.Range("D2:D51").NumberFormat = "General"
$ws.Range("D2:D51").Style = "Normal"
